# Insert a new weekly price record for "Black Amber" ciruela at row 174,
# pushing the existing "Crimsom fall" rows (old 174, 175) down to 175, 176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(174).Insert()

$ws.Range("A174").Value = 10
$ws.Range("B174").Value = "Vega Modelo de Temuco"
$ws.Range("C174").Value = "La Araucanía"
$ws.Range("D174").Value = 44595
$ws.Range("E174").Value = 9
$ws.Range("F174").Value = "Fruta"
$ws.Range("G174").Value = 100103
$ws.Range("H174").Value = "Frutos de hueso (carozo)"
$ws.Range("I174").Value = 100103002
$ws.Range("J174").Value = "Ciruela"
$ws.Range("K174").Value = "Black Amber"
$ws.Range("L174").Value = "Primera"
$ws.Range("M174").Value = 95
$ws.Range("N174").Value = 13000
$ws.Range("O174").Value = 13000
$ws.Range("P174").Value = 13000
$ws.Range("Q174").Value = "$/bandeja 18 kilos granel"
$ws.Range("R174").Value = "Región de O'Higgins"
$ws.Range("S174").Value = 722
$ws.Range("T174").Value = 18
